# Updated cryptos list on Thu Mar 23 06:51:24 UTC 2023 with GitHub Actions
# Applies the price / 1h-volume refresh + a PancakeSwap/Decentraland row swap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.694.45'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '1.754.74'
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'324.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.64%  '
$ws.Range("D6").Value = "'0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = "'0.4304"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.45%  '
$ws.Range("D8").Value = "'0.3638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.61%  '
$ws.Range("D9").Value = "'45.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").Value = "'0.07486"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = "'1.119"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.42%  '
$ws.Range("D12").Value = "'0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = "'21.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").Value = "'6.150"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.14%  '
$ws.Range("D15").Value = "'7.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.87%  '
$ws.Range("D16").Value = '1.748.16'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("D17").Value = "'0.00001068"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("D18").Value = "'87.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.53%  '
$ws.Range("D19").Value = "'0.06219"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.65%  '
$ws.Range("D20").Value = "'0.9992"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D22").Value = "'6.159"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.40%  '
$ws.Range("D23").Value = "'0.5271"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.69%  '
$ws.Range("D24").Value = '27.730.95'
$ws.Range("E24").Value = '  -1.95%  '
$ws.Range("D25").Value = "'11.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.13%  '
$ws.Range("D26").Value = "'2.325"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.39%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = "'152.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("D29").Value = "'2.366"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").Value = '1.948.13'
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").Value = "'1.219"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.68%  '
$ws.Range("D32").Value = "'127.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.33%  '
$ws.Range("D33").Value = "'5.723"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = "'0.09155"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.61%  '
$ws.Range("D35").Value = "'3.662"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.38%  '
$ws.Range("D36").Value = "'12.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.00%  '
$ws.Range("D37").Value = "'0.02311"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").Value = "'0.2153"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.72%  '
$ws.Range("D39").Value = "'5.111"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.17%  '
$ws.Range("D40").Value = "'0.6484"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").Value = "'0.06095"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.13%  '
$ws.Range("D42").Value = "'1.196"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.78%  '
$ws.Range("D44").Value = "'7.943"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.41%  '
$ws.Range("D45").Value = "'0.9990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").Value = "'13.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.14%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = "'3.751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.15%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.5942"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.27%  '
$ws.Range("D49").Value = "'126.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("D50").Value = "'1.973"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.33%  '
$ws.Range("D51").Value = "'0.06901"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.74%  '
